$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 is a new forecast record parsed from valid data, mirroring the
# shape of row 6 (same date/sector/passenger-type/airport/expected/pax
# figures) but the arrival time and flight number could not be parsed
# from the source, so those two cells are left blank (invalid data is
# ignored rather than causing the whole row to be skipped).
$src = $ws.Range("A6:J6")
$dst = $ws.Range("A7:J7")
$src.Copy($dst)

$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Move the active selection to reflect where the user would continue
# entering/reviewing data next.
$ws.Range("B8").Select()
